$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("H4").Value = 1684
$ws.Range("H7").Value = 25995
$ws.Range("I6").Value = 8966
$ws.Range("I7").Value = 26181
$ws.Range("J2").Value = 519
$ws.Range("J3").Value = 569
$ws.Range("J4").Value = 116
$ws.Range("J6").Value = 867
$ws.Range("J7").Value = 2113

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("J4").Value = 1
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("J3").Value = 27
$ws.Range("J6").Value = 27
$ws.Range("J7").Value = 75

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("J3").Value = 9
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("J3").Value = 34
$ws.Range("J6").Value = 23
$ws.Range("J7").Value = 73

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("J4").Value = 1
$ws.Range("J7").Value = 21

$ws = $wb.Worksheets.Item('New City')
$ws.Range("J2").Value = 10
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("H63").Value = 231
$ws.Range("H101").Value = 25995
$ws.Range("I12").Value = 66
$ws.Range("I101").Value = 26181
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 62
$ws.Range("J8").Value = 124
$ws.Range("J9").Value = 12
$ws.Range("J11").Value = 31
$ws.Range("J12").Value = 10
$ws.Range("J14").Value = 10
$ws.Range("J15").Value = 24
$ws.Range("J18").Value = 41
$ws.Range("J19").Value = 79
$ws.Range("J23").Value = 18
$ws.Range("J33").Value = 86
$ws.Range("J34").Value = 14
$ws.Range("J36").Value = 37
$ws.Range("J37").Value = 75
$ws.Range("J40").Value = 8
$ws.Range("J42").Value = 91
$ws.Range("J47").Value = 22
$ws.Range("J51").Value = 28
$ws.Range("J52").Value = 56
$ws.Range("J53").Value = 24
$ws.Range("J55").Value = 26
$ws.Range("J63").Value = 11
$ws.Range("J64").Value = 14
$ws.Range("J65").Value = 40
$ws.Range("J66").Value = 3
$ws.Range("J67").Value = 73
$ws.Range("J73").Value = 20
$ws.Range("J76").Value = 35
$ws.Range("J78").Value = 28
$ws.Range("J79").Value = 68
$ws.Range("J84").Value = 21
$ws.Range("J85").Value = 86
$ws.Range("J86").Value = 8
$ws.Range("J87").Value = 11
$ws.Range("J89").Value = 21
$ws.Range("J93").Value = 10
$ws.Range("J97").Value = 14
$ws.Range("J98").Value = 14
$ws.Range("J99").Value = 28
$ws.Range("J101").Value = 2113

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("J2").Value = 20
$ws.Range("J6").Value = 36
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("J6").Value = 37
$ws.Range("J7").Value = 79

$ws = $wb.Worksheets.Item('River North')
$ws.Range("J3").Value = 9
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("J2").Value = 17
$ws.Range("J3").Value = 28
$ws.Range("J4").Value = 8
$ws.Range("J7").Value = 86

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 17

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J2").Value = 14
$ws.Range("J6").Value = 61
$ws.Range("J7").Value = 91

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 26

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 18

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("J2").Value = 21
$ws.Range("J3").Value = 19
$ws.Range("J4").Value = 3
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 68

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("J2").Value = 2
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("J2").Value = 8
$ws.Range("J6").Value = 28
$ws.Range("J7").Value = 41

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J4").Value = 2
$ws.Range("J6").Value = 20
$ws.Range("J7").Value = 37

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("J2").Value = 1
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("J3").Value = 20
$ws.Range("J7").Value = 56

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("J6").Value = 9
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("J2").Value = 7
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range("J2").Value = 5
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("J6").Value = 1
$ws.Range("J7").Value = 3

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("J2").Value = 9
$ws.Range("J3").Value = 7
$ws.Range("J7").Value = 20

$ws = $wb.Worksheets.Item('West Town')
$ws.Range("J6").Value = 12
$ws.Range("J7").Value = 14

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J3").Value = 41
$ws.Range("J6").Value = 42
$ws.Range("J7").Value = 124

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 8

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J2").Value = 8
$ws.Range("J4").Value = 6
$ws.Range("J7").Value = 28

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("J6").Value = 15
$ws.Range("J7").Value = 24

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range("J2").Value = 3
$ws.Range("J7").Value = 8

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("J2").Value = 15
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 62

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("I6").Value = 36
$ws.Range("I7").Value = 66
$ws.Range("J6").Value = 8
$ws.Range("J7").Value = 10

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("J6").Value = 6
$ws.Range("J7").Value = 11
